$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function Get-ParaByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $para = $doc.Paragraphs.Item($i)
        $t = $para.Range.Text.TrimEnd([char]13)
        if ($t -eq $text) {
            return $para
        }
    }
    return $null
}

# ---------------------------------------------------------------
# 1) Heading "Data Structure Traveling with" + bookmark + " the Flow"
#    -> single run "Data Structure Traveling with the Flow", bookmark removed
# ---------------------------------------------------------------
$para = Get-ParaByText $d "Data Structure Traveling with the Flow"
$rng = $d.Range($para.Range.Start, $para.Range.End - 1)
$xml = "<w:p $wNs><w:pPr><w:pStyle w:val='Heading1'/><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>Data Structure Traveling with the Flow</w:t></w:r></w:p>"
$rng.InsertXML($xml) | Out-Null

# ---------------------------------------------------------------
# 2) "Alert Rule" paragraph: append " = [NIDS_Rule HIDS_Rule]" as 2nd run,
#    then insert 2 new paragraphs (NIDS_Rule=..., HIDS_Rule=) right after it.
# ---------------------------------------------------------------
$para = Get-ParaByText $d "Alert Rule"
$rng = $d.Range($para.Range.Start, $para.Range.End - 1)

$alertRuleXml = "<w:p $wNs>" +
    "<w:pPr><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr>" +
    "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>Alert Rule</w:t></w:r>" +
    "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> = [NIDS_Rule HIDS_Rule]</w:t></w:r>" +
    "</w:p>"

$nidsRuleXml = "<w:p $wNs>" +
    "<w:pPr><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr>" +
    "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'>NIDS_Rule = </w:t></w:r>" +
    "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>Action + P</w:t></w:r>" +
    "<w:r><w:t>rotocol</w:t></w:r>" +
    "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> + S</w:t></w:r>" +
    "<w:r><w:t>rc_ip</w:t></w:r>" +
    "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> + S</w:t></w:r>" +
    "<w:r><w:t>rc_port</w:t></w:r>" +
    "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> + D</w:t></w:r>" +
    "<w:r><w:t>irection</w:t></w:r>" +
    "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> + D</w:t></w:r>" +
    "<w:r><w:t>st_ip</w:t></w:r>" +
    "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> + D</w:t></w:r>" +
    "<w:r><w:t>st_port</w:t></w:r>" +
    "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> + (R</w:t></w:r>" +
    "<w:r><w:t>ule</w:t></w:r>" +
    "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>_</w:t></w:r>" +
    "<w:r><w:t>options)</w:t></w:r>" +
    "</w:p>"

$hidsRuleXml = "<w:p $wNs>" +
    "<w:pPr><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr>" +
    "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'>HIDS_Rule = </w:t></w:r>" +
    "</w:p>"

$rng.InsertXML($alertRuleXml + $nidsRuleXml + $hidsRuleXml) | Out-Null

# ---------------------------------------------------------------
# 3) "Comments" heading paragraph: add <w:lastRenderedPageBreak/> before the text run
# ---------------------------------------------------------------
$para = Get-ParaByText $d "Comments"
$rng = $d.Range($para.Range.Start, $para.Range.End - 1)
$xml = "<w:p $wNs><w:pPr><w:pStyle w:val='Heading1'/><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:lastRenderedPageBreak/><w:t>Comments</w:t></w:r></w:p>"
$rng.InsertXML($xml) | Out-Null

# ---------------------------------------------------------------
# 4) After "The volume/time may not be correct" paragraph, insert new final
#    paragraph "HIDS_Rule data structure" + the _GoBack bookmark.
# ---------------------------------------------------------------
$para = Get-ParaByText $d "The volume/time may not be correct"
$rng = $d.Range($para.Range.Start, $para.Range.End - 1)
$xml = "<w:p $wNs>" +
    "<w:pPr><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr>" +
    "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>HIDS_Rule data structure</w:t></w:r>" +
    "<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/>" +
    "</w:p>"
$rng.InsertXML($xml) | Out-Null

# ---------------------------------------------------------------
# 5) Table cell: "Thang" / " Duong Chi" (2 runs + proofErr) -> single run
#    "Thang Duong Chi". Done LAST: the runtime's Paragraphs index cache
#    gets invalidated after mutating a table cell via InsertXML, so no
#    further Paragraphs-index-based lookups should happen after this.
# ---------------------------------------------------------------
$tbl = $d.Tables.Item(1)
$cell = $tbl.Cell(2, 1)
$cellPara = $cell.Range.Paragraphs.Item(1)
$cellRng = $d.Range($cellPara.Range.Start, $cellPara.Range.End - 1)
$xml = "<w:p $wNs><w:pPr><w:jc w:val='center'/></w:pPr><w:r><w:t>Thang Duong Chi</w:t></w:r></w:p>"
$cellRng.InsertXML($xml) | Out-Null

Write-Host "done"
